$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.358.42'
$ws.Range('E2').Value = '  +0.62%  '

$ws.Range('D3').Value = '1.878.19'
$ws.Range('E3').Value = '  +0.70%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.001'
$ws.Range('E4').Value = '  +0.12%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '244.25'
$ws.Range('E5').Value = '  +4.25%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.001'
$ws.Range('E6').Value = '  +0.06%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4764'
$ws.Range('E7').Value = '  +1.56%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2880'
$ws.Range('E8').Value = '  +1.22%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '42.65'
$ws.Range('E9').Value = '  +2.86%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.06527'
$ws.Range('E10').Value = '  -0.48%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '21.35'
$ws.Range('E11').Value = '  +0.15%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07766'
$ws.Range('E12').Value = '  +0.09%  '

$ws.Range('B13').Value = 'Polygon'
$ws.Range('C13').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.7395'
$ws.Range('E13').Value = '  +7.24%  '

$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').Value = '1.879.66'
$ws.Range('E14').Value = '  +1.19%  '

$ws.Range('B15').Value = 'Litecoin'
$ws.Range('C15').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '96.39'
$ws.Range('E15').Value = '  +0.48%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '5.136'
$ws.Range('E16').Value = '  +0.81%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '276.28'
$ws.Range('E17').Value = '  +4.22%  '

$ws.Range('D18').Value = '30.354.77'

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.39'
$ws.Range('E19').Value = '  -1.78%  '

$ws.Range('B20').Value = 'ShibaInu'
$ws.Range('C20').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.000007541'
$ws.Range('E20').Value = '  -1.89%  '

$ws.Range('B21').Value = 'Dai'
$ws.Range('C21').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.001'
$ws.Range('E21').Value = '  +0.06%  '

$ws.Range('D22').Value = '2.127.89'
$ws.Range('E22').Value = '  +0.87%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.001'
$ws.Range('E23').Value = '  +0.17%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.229'
$ws.Range('E24').Value = '  -0.06%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '6.171'
$ws.Range('E25').Value = '  +0.23%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.235'
$ws.Range('E26').Value = '  -2.31%  '

$ws.Range('E27').Value = '  -1.08%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.93'
$ws.Range('E28').Value = '  +1.44%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.955'
$ws.Range('E29').Value = '  +1.19%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.382'
$ws.Range('E30').Value = '  +0.95%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.09958'
$ws.Range('E31').Value = '  +0.62%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.513'
$ws.Range('E32').Value = '  +3.86%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.312'
$ws.Range('E33').Value = '  -1.18%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.086'
$ws.Range('E34').Value = '  +1.35%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.04751'
$ws.Range('E35').Value = '  +0.80%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.122'
$ws.Range('E36').Value = '  -0.46%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.6955'
$ws.Range('E37').Value = '  -0.34%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.719'
$ws.Range('E38').Value = '  +0.12%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01852'
$ws.Range('E39').Value = '  -0.03%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.750'
$ws.Range('E40').Value = '  -0.70%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.300'
$ws.Range('E41').Value = '  +0.72%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.8427'
$ws.Range('E42').Value = '  +1.15%  '

$ws.Range('B43').Value = 'TheSandbox'
$ws.Range('C43').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.4175'
$ws.Range('E43').Value = '  +1.14%  '

$ws.Range('B44').Value = 'RenderToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.909'
$ws.Range('E44').Value = '  -0.99%  '

$ws.Range('B45').Value = 'Aave'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '69.35'
$ws.Range('E45').Value = '  -4.01%  '

$ws.Range('B46').Value = 'PaxDollar'
$ws.Range('C46').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.000'
$ws.Range('E46').Value = '  +0.00%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '101.73'
$ws.Range('E47').Value = '  -0.85%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '9.258'
$ws.Range('E48').Value = '  +2.00%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.098'
$ws.Range('E49').Value = '  +0.31%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '35.18'
$ws.Range('E50').Value = '  +1.90%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '910.63'
$ws.Range('E51').Value = '  -6.25%  '
